$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.251.48"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "'2.974.90"
$ws.Range("E3").Value = "  +1.79%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'383.23"
$ws.Range("E5").Value = "  +1.58%  "

$ws.Range("D6").Value = "'102.00"
$ws.Range("E6").Value = "  -0.91%  "

$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "'0.591"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").Value = "'36.60"
$ws.Range("E10").Value = "  -1.20%  "

$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").Value = "'0.0840"
$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("D13").Value = "'3.451.75"
$ws.Range("E13").Value = "  +1.94%  "

$ws.Range("D14").Value = "'18.10"
$ws.Range("E14").Value = "  -1.11%  "

$ws.Range("D15").Value = "'7.46"
$ws.Range("E15").Value = "  +1.60%  "

$ws.Range("D16").Value = "'2.982.11"
$ws.Range("E16").Value = "  +2.18%  "

$ws.Range("D17").Value = "'0.993"
$ws.Range("E17").Value = "  +7.06%  "

$ws.Range("D18").Value = "'51.234.68"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("D19").Value = "'3.23"
$ws.Range("E19").Value = "  -5.63%  "

$ws.Range("D20").Value = "'7.35"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").Value = "'12.70"
$ws.Range("E21").Value = "  -1.75%  "

$ws.Range("D22").Value = "'0.0₃0955"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("D23").Value = "'68.77"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").Value = "'261.65"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("D25").Value = "'2.89"
$ws.Range("E25").Value = "  +4.66%  "

$ws.Range("D26").Value = "'8.16"
$ws.Range("E26").Value = "  +13.37%  "

$ws.Range("D27").Value = "'7.57"
$ws.Range("E27").Value = "  +10.92%  "

$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").Value = "'0.113"
$ws.Range("E28").Value = "  +12.26%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.167"
$ws.Range("E29").Value = "  -1.70%  "

$ws.Range("D30").Value = "'4.10"
$ws.Range("E30").Value = "  -0.59%  "

$ws.Range("E31").Value = "  -0.01%  "

$ws.Range("D32").Value = "'25.77"
$ws.Range("E32").Value = "  +0.47%  "

$ws.Range("D33").Value = "'9.83"
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("D34").Value = "'34.21"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("D35").Value = "'50.90"
$ws.Range("E35").Value = "  -0.81%  "

$ws.Range("E36").Value = "  -2.66%  "

$ws.Range("D37").Value = "'0.0445"
$ws.Range("E37").Value = "  +5.75%  "

$ws.Range("D39").Value = "'2.98"
$ws.Range("E39").Value = "  -0.84%  "

$ws.Range("D40").Value = "'17.00"
$ws.Range("E40").Value = "  +0.60%  "

$ws.Range("D41").Value = "'2.57"
$ws.Range("E41").Value = "  +0.89%  "

$ws.Range("E42").Value = "  +1.22%  "

$ws.Range("D43").Value = "'1.79"
$ws.Range("E43").Value = "  -1.06%  "

$ws.Range("D44").Value = "'122.38"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").Value = "'21.36"
$ws.Range("E45").Value = "  -0.76%  "

$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("D47").Value = "'0.275"
$ws.Range("E47").Value = "  +1.49%  "

$ws.Range("E48").Value = "  +1.84%  "

$ws.Range("D49").Value = "'2.024.64"
$ws.Range("E49").Value = "  -0.04%  "

$ws.Range("D50").Value = "'3.23"
$ws.Range("E50").Value = "  +2.62%  "

$ws.Range("D51").Value = "'0.0333"
$ws.Range("E51").Value = "  +5.47%  "
